# Auto-generated edit script applying scheduled-runner profit recalculations
# to the Aegis_Profits workbook across ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 833.2889
$ws.Range("I17").Value = 297
$ws.Range("J17").Value = 845.4773
$ws.Range("K17").Value = 891
$ws.Range("L17").Value = 2536.4319
$ws.Range("M17").Value = -723
$ws.Range("N17").Value = -2872.4319

$ws.Range("H19").Value = 2062.9092
$ws.Range("I19").Value = 1925
$ws.Range("J19").Value = 2141.7144
$ws.Range("K19").Value = 1925
$ws.Range("L19").Value = 2141.7144
$ws.Range("M19").Value = -1750
$ws.Range("N19").Value = -2491.7144

$ws.Range("H46").Value = 929.4483
$ws.Range("J46").Value = 929.4483
$ws.Range("L46").Value = 2788.3449
$ws.Range("N46").Value = -3026.3449

$ws.Range("H60").Value = 929.4483
$ws.Range("J60").Value = 929.4483
$ws.Range("L60").Value = 2788.3449
$ws.Range("N60").Value = -3756.3449

$ws.Range("H129").Value = 2158.9312
$ws.Range("J129").Value = 1105.2361
$ws.Range("L129").Value = 3315.7083
$ws.Range("N129").Value = -13315.7083

$ws.Range("H138").Value = 1860.3405
$ws.Range("J138").Value = 3439.1428
$ws.Range("L138").Value = 10317.4284
$ws.Range("N138").Value = -20597.4284

$ws.Range("H140").Value = 53789.5
$ws.Range("J140").Value = 53789.5
$ws.Range("L140").Value = 53789.5
$ws.Range("N140").Value = -64149.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20648.031
$ws.Range("I32").Value = 4629.1787
$ws.Range("K32").Value = 4629.1787
$ws.Range("M32").Value = -4342.1787

$ws.Range("H61").Value = 1447.8966
$ws.Range("I61").Value = 1276.8182
$ws.Range("J61").Value = 1985.5714
$ws.Range("K61").Value = 1276.8182
$ws.Range("L61").Value = 1985.5714
$ws.Range("M61").Value = -1064.8182
$ws.Range("N61").Value = -2409.5714

$ws.Range("H102").Value = 73819.21000000001
$ws.Range("I102").Value = 113821
$ws.Range("J102").Value = 1816
$ws.Range("K102").Value = 113821
$ws.Range("L102").Value = 1816
$ws.Range("M102").Value = -112199
$ws.Range("N102").Value = -5060

$ws.Range("H136").Value = 1447.8966
$ws.Range("I136").Value = 1276.8182
$ws.Range("J136").Value = 1985.5714
$ws.Range("K136").Value = 3830.4546
$ws.Range("L136").Value = 5956.7142
$ws.Range("M136").Value = -1280.4546
$ws.Range("N136").Value = -11056.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 32801.906
$ws.Range("I20").Value = 43306.043
$ws.Range("J20").Value = 1289.5
$ws.Range("K20").Value = 43306.043
$ws.Range("L20").Value = 1289.5
$ws.Range("M20").Value = -43059.043
$ws.Range("N20").Value = -1783.5

$ws.Range("H86").Value = 58410.55
$ws.Range("I86").Value = 76813.60000000001
$ws.Range("J86").Value = 3201.4
$ws.Range("K86").Value = 76813.60000000001
$ws.Range("L86").Value = 3201.4
$ws.Range("M86").Value = -75690.60000000001
$ws.Range("N86").Value = -5447.4

$ws.Range("H89").Value = 58410.55
$ws.Range("I89").Value = 76813.60000000001
$ws.Range("J89").Value = 3201.4
$ws.Range("K89").Value = 384068
$ws.Range("L89").Value = 16007
$ws.Range("M89").Value = -378452
$ws.Range("N89").Value = -27239

$ws.Range("H107").Value = 166734080
$ws.Range("I107").Value = 166734080
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 166734080
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -166732160

$ws.Range("H133").Value = 359999
$ws.Range("J133").Value = 359999
$ws.Range("L133").Value = 359999
$ws.Range("N133").Value = -370119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 330.8
$ws.Range("I22").Value = 330.25
$ws.Range("K22").Value = 330.25
$ws.Range("M22").Value = 19.75

$ws.Range("H31").Value = 43344.94
$ws.Range("I31").Value = 2431.8572
$ws.Range("J31").Value = 70620.336
$ws.Range("K31").Value = 2431.8572
$ws.Range("L31").Value = 70620.336
$ws.Range("M31").Value = -2136.8572
$ws.Range("N31").Value = -71210.336

$ws.Range("H34").Value = 43344.94
$ws.Range("I34").Value = 2431.8572
$ws.Range("J34").Value = 70620.336
$ws.Range("K34").Value = 2431.8572
$ws.Range("L34").Value = 70620.336
$ws.Range("M34").Value = -2229.8572
$ws.Range("N34").Value = -71024.336

$ws.Range("H132").Value = 3836.147
$ws.Range("I132").Value = 3681.3462
$ws.Range("J132").Value = 4339.25
$ws.Range("K132").Value = 11044.0386
$ws.Range("L132").Value = 13017.75
$ws.Range("M132").Value = -8514.0386
$ws.Range("N132").Value = -18077.75

$ws.Range("H134").Value = 1244.4
$ws.Range("I134").Value = 1235.9412
$ws.Range("J134").Value = 1262.375
$ws.Range("K134").Value = 3707.8236
$ws.Range("L134").Value = 3787.125
$ws.Range("M134").Value = -1172.8236
$ws.Range("N134").Value = -8857.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 7094.1763
$ws.Range("I117").Value = 321.75
$ws.Range("J117").Value = 9178
$ws.Range("K117").Value = 965.25
$ws.Range("L117").Value = 27534
$ws.Range("M117").Value = 2476.75
$ws.Range("N117").Value = -34418

$ws.Range("H131").Value = 1404.2892
$ws.Range("I131").Value = 1451
$ws.Range("J131").Value = 1401.9241
$ws.Range("K131").Value = 4353
$ws.Range("L131").Value = 4205.7723
$ws.Range("M131").Value = 687
$ws.Range("N131").Value = -14285.7723

$ws.Range("H138").Value = 10511.667
$ws.Range("I138").Value = 13015.556
$ws.Range("K138").Value = 39046.66800000001
$ws.Range("M138").Value = -33906.66800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 68857.23
$ws.Range("I70").Value = 94944.55
$ws.Range("J70").Value = 5088.222
$ws.Range("K70").Value = 94944.55
$ws.Range("L70").Value = 5088.222
$ws.Range("M70").Value = -94674.55
$ws.Range("N70").Value = -5628.222

$ws.Range("H73").Value = 68857.23
$ws.Range("I73").Value = 94944.55
$ws.Range("J73").Value = 5088.222
$ws.Range("K73").Value = 94944.55
$ws.Range("L73").Value = 5088.222
$ws.Range("M73").Value = -94008.55
$ws.Range("N73").Value = -6960.222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3830
$ws.Range("I7").Value = 2200
$ws.Range("J7").Value = 4528.5713
$ws.Range("K7").Value = 2200
$ws.Range("L7").Value = 4528.5713
$ws.Range("M7").Value = -2088
$ws.Range("N7").Value = -4752.5713

$ws.Range("H46").Value = 723389.1
$ws.Range("I46").Value = 314.66666
$ws.Range("J46").Value = 1265695
$ws.Range("K46").Value = 314.66666
$ws.Range("L46").Value = 1265695
$ws.Range("M46").Value = -126.66666
$ws.Range("N46").Value = -1266071

$ws.Range("H55").Value = 308054.38
$ws.Range("I55").Value = 1033324.2
$ws.Range("K55").Value = 1033324.2
$ws.Range("M55").Value = -1033151.2

$ws.Range("H61").Value = 1894.5834
$ws.Range("I61").Value = 1886.4286
$ws.Range("J61").Value = 1906
$ws.Range("K61").Value = 1886.4286
$ws.Range("L61").Value = 1906
$ws.Range("M61").Value = -1684.4286
$ws.Range("N61").Value = -2310

$ws.Range("H113").Value = 1894.5834
$ws.Range("I113").Value = 1886.4286
$ws.Range("J113").Value = 1906
$ws.Range("K113").Value = 1886.4286
$ws.Range("L113").Value = 1906
$ws.Range("M113").Value = 283.5714
$ws.Range("N113").Value = -6246

$ws.Range("H126").Value = 3830
$ws.Range("I126").Value = 2200
$ws.Range("J126").Value = 4528.5713
$ws.Range("K126").Value = 6600
$ws.Range("L126").Value = 13585.7139
$ws.Range("M126").Value = -4130
$ws.Range("N126").Value = -18525.7139

$ws.Range("H132").Value = 4004.5
$ws.Range("I132").Value = 4787.125
$ws.Range("J132").Value = 1917.5
$ws.Range("K132").Value = 14361.375
$ws.Range("L132").Value = 5752.5
$ws.Range("M132").Value = -11831.375
$ws.Range("N132").Value = -10812.5

$ws.Range("H136").Value = 1033.6296
$ws.Range("I136").Value = 855.449
$ws.Range("J136").Value = 2779.8
$ws.Range("K136").Value = 2566.347
$ws.Range("L136").Value = 8339.400000000001
$ws.Range("M136").Value = -16.34699999999975
$ws.Range("N136").Value = -13439.4
